$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 0.1510678514225969
$ws.Range("M2").Value = 2.425633666666667
$ws.Range("N2").Value = 7.276901000000001
$ws.Range("O2").Value = 0.0662600404061536
$ws.Range("P2").Value = 0.06626004040615362
$ws.Range("Q2").Value = 0.3514104432801111
$ws.Range("R2").Value = 3.162693989521
$ws.Range("S2").Value = 0.01000976193933208
$ws.Range("T2").Value = 0.01000976193933208

# Row 3
$ws.Range("J3").Value = 0.1510678514225969
$ws.Range("O3").Value = 0.4234968256437875
$ws.Range("P3").Value = 0.4234968256437876
$ws.Range("S3").Value = 0.06397675553429712
$ws.Range("T3").Value = 0.06397675553429714

# Row 4
$ws.Range("J4").Value = 0.1510678514225969
$ws.Range("M4").Value = 18.67887366666667
$ws.Range("N4").Value = 56.03662100000001
$ws.Range("O4").Value = 0.5102431339500588
$ws.Range("P4").Value = 0.5102431339500588
$ws.Range("Q4").Value = 2.706076917293445
$ws.Range("R4").Value = 24.35469225564101
$ws.Range("S4").Value = 0.0770813339489677
$ws.Range("T4").Value = 0.07708133394896771

# Row 5
$ws.Range("G5").Value = 0.8141236666666667
$ws.Range("M5").Value = 2.425633666666667
$ws.Range("N5").Value = 7.276901000000001
$ws.Range("O5").Value = 0.0662600404061536
$ws.Range("P5").Value = 0.06626004040615362
$ws.Range("Q5").Value = 1.974765774696778
$ws.Range("R5").Value = 17.772891972271
$ws.Range("S5").Value = 0.05625027846682152
$ws.Range("T5").Value = 0.05625027846682153

# Row 6
$ws.Range("G6").Value = 0.8141236666666667
$ws.Range("O6").Value = 0.4234968256437875
$ws.Range("P6").Value = 0.4234968256437876
$ws.Range("S6").Value = 0.3595200701094904
$ws.Range("T6").Value = 0.3595200701094904

# Row 7
$ws.Range("G7").Value = 0.8141236666666667
$ws.Range("M7").Value = 18.67887366666667
$ws.Range("N7").Value = 56.03662100000001
$ws.Range("O7").Value = 0.5102431339500588
$ws.Range("P7").Value = 0.5102431339500588
$ws.Range("Q7").Value = 15.20691311871012
$ws.Range("S7").Value = 0.4331618000010911
$ws.Range("T7").Value = 0.4331618000010911
